# Update "想去人数" (want-to-go count) values that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 5519
$wsExhibit.Range("F7").Value = 33
$wsExhibit.Range("F8").Value = 361

# Sheet "演出" (performances)
$wsPerform = $wb.Worksheets.Item("演出")
$wsPerform.Range("F2").Value = 48

# Sheet "全部类型" (all types, combined view)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 5519
$wsAll.Range("F7").Value = 33
$wsAll.Range("F8").Value = 48
$wsAll.Range("F9").Value = 361
